$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999207303236,
    0.9951792243656219,
    0.97310449354629558,
    0.9648505055369565,
    0.95706098909194526,
    0.94117971484604479,
    0.93778405841862233,
    0.93415576778690967,
    0.9315478592080213,
    0.92988237057145118,
    0.92977994558757648,
    0.92984139938621035,
    0.93167504139234103,
    0.93318985971473745,
    0.93519206289179646,
    0.93830839387214393,
    0.93460081222559954,
    0.93349193436346145,
    0.99412602342080836,
    0.98700922226829091,
    0.98561078332973218,
    0.98434630123655187,
    0.96606894022155643,
    0.95304728973569253,
    0.94659017544450219,
    0.93921089014612913,
    0.93540221122774403,
    0.92022428485011454,
    0.90960052813372327,
    0.90478769180199303,
    0.90339914367534091,
    0.90406033610299819,
    0.90354034591190824
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

